$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.014.98"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.908.28"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.98"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4829"
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3806"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07361"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9340"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.81"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07787"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.920.16"
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.495"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.650"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.83"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008880"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "28.047.03"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.73"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.156"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.178.69"
$ws.Range("E23").Value = "  +3.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.90"
$ws.Range("E24").Value = "  +2.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.27"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("E26").Value = "  -1.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.51"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.120"
$ws.Range("E28").Value = "  +5.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.17"
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.986"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08944"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.287"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.255"
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7722"
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.665"
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.599"
$ws.Range("E36").Value = "  -3.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02054"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5521"
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05296"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.994"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.003"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.501"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "109.92"
$ws.Range("E45").Value = "  +6.88%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.74"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4824"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.17"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06078"
$ws.Range("E51").Value = "  +0.04%  "
